$d = $word.ActiveDocument
$br = [char]11

# ---------------------------------------------------------------------------
# 1) "Ativação: 15/07/2024" -> "Ativação: Semestral"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 15/07/2024", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: Semestral", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Objetivos paragraph: "Introduzir aos estudantes..." -> "01/01/2025"
#    (must run before step 10 re-introduces the same sentence elsewhere)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Introduzir aos estudantes os princípios e a metodologia da pesquisa científica.", $true, $false, $false, $false, $false, $true, 1, $false, "01/01/2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Programa paragraph (old typo "pesquia") -> Bibliografia-style reference list
#    (must run before step 10 inserts the corrected "pesquisa" wording)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Organização e o formalismo do desenvolvimento do trabalho científico ou projeto de engenharia. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica. Desenvolvimento de um tema de pesquia ou projeto de engenharia, com o formato de um trabalho de iniciação científica, sob a orientação de um professor ou pesquisador autorizado pela Comissão de Curso. Entrega e apresentação de documento técnico no final da disciplina.", $true, $false, $false, $false, $false, $true, 1, $false, "ASTI VERA, A. Metodologia da pesquisa científica. Porto Alegre: Ed. Globo, 1973. BARRAS, R. Os cientistas precisam escrever: guia de redação para cientistas, engenheiros e estudantes. São Paulo: TAQ/EDUSP, 1979. CERVO, A. L.; BERVIAN, P. A. Metodologia científica. São Paulo: Mc-Graw-Hill do Brasil, 1973. ANDRADE, M. M. Introdução à Metodologia do Trabalho Científico São Paulo: Atlas, 2005.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Programa resumido paragraph: "Iniciação a um projeto..." -> recuperação note
#    (must run before step 10 re-introduces the same sentence elsewhere)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Iniciação a um projeto de pesquisa sob orientação de um professor.", $true, $false, $false, $false, $false, $true, 1, $false, "Devido às características práticas da disciplina, não será oferecida recuperação", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Italic EN summary: "Initiation into..." -> "Initiation to..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Initiation into a research project under the guidance of a professor.", $true, $false, $false, $false, $false, $true, 1, $false, "Initiation to a research project under the guidance of a professor.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Italic EN Programa: "engineering project...research project document" -> "engineering design...technical document"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Organization and formalism of the development of scientific work or engineering project. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research or engineering project topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project document at the end of the course.", $true, $false, $false, $false, $false, $true, 1, $false, "Organization and formalism of the development of scientific work or engineering design. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of a research topic or engineering project, in the format of a scientific initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of technical document at the end of the course.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Avaliação / Método value: "Aulas expositivas..." -> "Estudantes de ensino infantil, fundamental ou médio."
#    (must run before step 10 re-introduces the same sentence elsewhere)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e/ou engenharia e elaboração de projeto de pesquisa e/ou engenharia. Visitas técnicas em institutos ou empresas da área científica ou de engenharia.", $true, $false, $false, $false, $false, $true, 1, $false, "Estudantes de ensino infantil, fundamental ou médio.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Avaliação / Critério value: "Nota de avaliação do projeto e demais documentos." -> new text
#    (must run before step 10 re-introduces the same sentence elsewhere)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Nota de avaliação do projeto e demais documentos.", $true, $false, $false, $false, $false, $true, 1, $false, "Para os estudantes: despertar interesse na engenharia.Para a formação dos discentes: Desenvolver conceitos de engenharia com aplicações profissionais", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Avaliação / Norma de recuperação value: "Devido às características do curso..." -> long bullet text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Devido às características do curso, não será oferecida recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "- Identificação das necessidades do grupo social: pesquisas, entrevistas e observações para entender as necessidades, desafios e preferências dos estudantes.- Definição de objetivos e requisitos do projeto para que as soluções desenvolvidas devem atender: identificar funcionalidades, restrições de orçamento e cronograma, e quaisquer outras considerações importantes.- Pesquisa e desenvolvimento projetos relacionados à engenharia: criação de protótipos, desenvolvimento de software, fabricação de dispositivos e apresentação de aplicações para garantir que haja disseminação do conhecimento sobre a profissão engenharia.- Avaliação: feedback recebido quanto ao conhecimento sobre o tema.- Implementação e distribuição: Visita e apresentações em escolas de ensino infantil, fundamental ou médio.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10) Bibliografia paragraph: collapse the 4-reference run/break list into a
#     single new sentence (whole-paragraph text replace collapses to one run).
# ---------------------------------------------------------------------------
$pBiblio = $d.Paragraphs.Item(19)
$pBiblio.Range.Text = "Os indicadores serão obtidos por questionário de avaliação pelos usuários quanto aos seguintes quesitos: conhecimento adquirido e satisfação do usuário nas apresentações e formas de divulgação."

# ---------------------------------------------------------------------------
# 11) Docente(s) bullet-list paragraph: splice in four new lines before
#     "519033" and append three more after "Katia..." — this reuses text that
#     has already been moved out of its old homes above, so no ambiguous
#     matches remain at this point.
# ---------------------------------------------------------------------------

# 11a) Insert "Introduzir aos estudantes..." + line break right before "519033"
$paraObj = $d.Paragraphs.Item(9)
$rStart = $paraObj.Range
$insStart = $d.Range($rStart.Start, $rStart.Start)
$txtIntro = "Introduzir aos estudantes os princípios e a metodologia da pesquisa científica." + $br
$insStart.InsertBefore($txtIntro)

# 11b) Append a line break to the end of "5817692 - Katia..." run (merges into
#      that same run, matching the reference structure).
$d.Content.Find.Execute("5817692 - Katia Cristiane Gandolpho Candioto", $true, $false, $false, $false, $false, $true, 1, $false, "5817692 - Katia Cristiane Gandolpho Candioto^l", 2) | Out-Null

# 11c) Append the remaining new lines as their own runs, one at a time, by
#      inserting right before the paragraph mark each time.
$paraObj = $d.Paragraphs.Item(9)
$pos1 = $paraObj.Range.End - 1
$insR1 = $d.Range($pos1, $pos1)
$txt2 = "Iniciação a um projeto de pesquisa sob orientação de um professor." + $br
$insR1.InsertAfter($txt2)

$paraObj = $d.Paragraphs.Item(9)
$pos2 = $paraObj.Range.End - 1
$insR2 = $d.Range($pos2, $pos2)
$txt3 = "Organização e o formalismo do desenvolvimento do trabalho científico ou projeto de engenharia. Técnicas de redação científica, uso de ferramentas de busca, referências bibliográficas e estruturas formais de divulgação científica. Desenvolvimento de um tema de pesquisa ou projeto de engenharia, com o formato de um trabalho de iniciação científica, sob a orientação de um professor ou pesquisador autorizado pela Comissão de Curso. Entrega e apresentação de documento técnico no final da disciplina." + $br
$insR2.InsertAfter($txt3)

$paraObj = $d.Paragraphs.Item(9)
$pos3 = $paraObj.Range.End - 1
$insR3 = $d.Range($pos3, $pos3)
$txt4 = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e/ou engenharia e elaboração de projeto de pesquisa e/ou engenharia. Visitas técnicas em institutos ou empresas da área científica ou de engenharia." + $br
$insR3.InsertAfter($txt4)

$paraObj = $d.Paragraphs.Item(9)
$pos4 = $paraObj.Range.End - 1
$insR4 = $d.Range($pos4, $pos4)
$txt5 = "Nota de avaliação do projeto e demais documentos."
$insR4.InsertAfter($txt5)

Write-Output "Done"
